$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Delete row 13 (the old "5840535 - Messias Borges Silva" answer row that
#    sat, label-less, right under "Docentes responsaveis:"). This shifts
#    rows 14-24 up to 13-23 and shrinks the used range to A1:C23, matching
#    every row-height / label-position change in the diff for free.
$ws.Rows.Item(13).Delete()

# 2) After the shift, a handful of answer cells (column B and the mirrored
#    column C) have different text than what simply shifted into place, so
#    fix those up individually.

# Objetivos: long description -> short "matricula - nome" answer
$ws.Range("B10").Value = "5840535 - Messias Borges Silva"
$ws.Range("C10").Value = "5840535 - Messias Borges Silva"

# Programa resumido: long syllabus text -> "Semestral"
$ws.Range("B13").Value = "Semestral"
$ws.Range("C13").Value = "Semestral"

# Programa: long syllabus text -> activation date. Copy the existing
# "01/01/1996" text cell (row 8) so it stays a real text value instead of
# Excel auto-converting the literal string into a date serial number.
$ws.Range("B8").Copy($ws.Range("B15"))
$ws.Range("C8").Copy($ws.Range("C15"))

# Metodo: (previously blank) -> "matricula - nome" answer
$ws.Range("B18").Value = "5840535 - Messias Borges Silva"
$ws.Range("C18").Value = "5840535 - Messias Borges Silva"

# Criterio: long evaluation text -> "2 provas escritas"
$ws.Range("B19").Value = "2 provas escritas"
$ws.Range("C19").Value = "2 provas escritas"

# Norma de recuperacao: recovery text -> long evaluation text
$ws.Range("B20").Value = "serão avaliados os conteúdos discutidos em sala e constantes da ementa do curso.A média da disciplina será a média aritmética das duas provas."
$ws.Range("C20").Value = "serão avaliados os conteúdos discutidos em sala e constantes da ementa do curso.A média da disciplina será a média aritmética das duas provas."

# Bibliografia: bibliography text -> recovery text
$ws.Range("B21").Value = "uma provas escrita com conteúdo de todo o semestre"
$ws.Range("C21").Value = "uma provas escrita com conteúdo de todo o semestre"
